# Apply roster update to the active worksheet.
# Rows 2-19 (columns A:C) are rewritten with the updated player/position/team
# data: the roster was reshuffled, "P.J. Washington" (PF, Dallas Mavericks)
# was dropped, and "Jrue Holiday" (PG,SG, Boston Celtics) was added.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Jalen Brunson",      "PG",       "New York Knicks"),
    @("Trae Young",         "PG",       "Atlanta Hawks"),
    @("Draymond Green",     "PF,C",     "Golden State Warriors"),
    @("Alperen Sengün",     "C",        "Houston Rockets"),
    @("Dereck Lively II",   "C",        "Dallas Mavericks"),
    @("Walker Kessler",     "C",        "Utah Jazz"),
    @("Yves Missi",         "C",        "New Orleans Pelicans"),
    @("Norman Powell",      "SG,SF",    "LA Clippers"),
    @("Jeremy Sochan",      "SF,PF",    "San Antonio Spurs"),
    @("Shaedon Sharpe",     "SG,SF",    "Portland Trail Blazers"),
    @("Devin Vassell",      "SG,SF",    "San Antonio Spurs"),
    @("Coby White",         "PG,SG",    "Chicago Bulls"),
    @("Devin Booker",       "PG,SG",    "Phoenix Suns"),
    @("LeBron James",       "SF,PF",    "Los Angeles Lakers"),
    @("Jrue Holiday",       "PG,SG",    "Boston Celtics"),
    @("Immanuel Quickley",  "PG,SG",    "Toronto Raptors"),
    @("Kawhi Leonard",      "SG,SF,PF", "LA Clippers"),
    @("Desmond Bane",       "SG,SF",    "Memphis Grizzlies")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
